$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update the volumenMayoristaReal (D) values that changed in this revision
# ---------------------------------------------------------------------------
$ws.Range("D6").Value = 5033240
$ws.Range("D9").Value = 5013010
$ws.Range("D14").Value = 1078000
$ws.Range("D15").Value = 1789750
$ws.Range("D20").Value = 2533993
$ws.Range("D21").Value = 369330

# D8 and D22 also lose their old "#,##0" number formatting in this revision
$ws.Range("D8").Value = 18756600
$ws.Range("D8").ClearFormats()
$ws.Range("D22").Value = 3656760
$ws.Range("D22").ClearFormats()

# ---------------------------------------------------------------------------
# 2. New header labels for the added columns G, H, I, K
# ---------------------------------------------------------------------------
$ws.Range("G1").Value = "diferencia2"
$ws.Range("H1").Value = "diferencia3"
$ws.Range("I1").Value = "diferencia4"
$ws.Range("K1").Value = "diferencia4"

# ---------------------------------------------------------------------------
# 3. New formula columns G, H, I (row 2 entered explicitly, rows 3-22 filled
#    as one range so the engine records them as a shared formula, matching
#    the authoring pattern already used by columns E/F)
# ---------------------------------------------------------------------------
$ws.Range("G2").Formula = "=SUM(D2,-B2)/D2"
$ws.Range("G3:G22").Formula = "=SUM(D3,-B3)/D3"

$ws.Range("H2").Formula = "=SUM((B2/D2)*100,-100)"
$ws.Range("H3:H22").Formula = "=SUM((B3/D3)*100,-100)"

$ws.Range("I2").Formula = "=SUM((B2/D2),-1)"
$ws.Range("I3:I22").Formula = "=SUM((B3/D3),-1)"

# formulas referencing D pick up D's own (non-default) number format on a
# couple of rows (D10/D13 still use the "#,##0" style) - strip that back off
# so G/H/I stay on the default/general format, matching the source edit
$ws.Range("G2:G22").ClearFormats()
$ws.Range("H2:H22").ClearFormats()
$ws.Range("I2:I22").ClearFormats()

# ---------------------------------------------------------------------------
# 4. Column K holds the same numbers as column I but pasted as plain values
# ---------------------------------------------------------------------------
for ($row = 2; $row -le 22; $row++) {
    $ws.Range("K$row").Value = $ws.Range("I$row").Value()
}

# ---------------------------------------------------------------------------
# 5. Number formatting: "0.0" applied to the diferencia (E) and K columns
# ---------------------------------------------------------------------------
$ws.Range("E1:E22").NumberFormat = "0.0"
$ws.Range("K1:K22").NumberFormat = "0.0"

# ---------------------------------------------------------------------------
# 6. Column widths
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 26
$ws.Columns.Item(4).ColumnWidth = 23.67
$ws.Columns.Item(6).ColumnWidth = 15.33
$ws.Columns.Item(7).ColumnWidth = 15

# ---------------------------------------------------------------------------
# 7. Page orientation + selection
# ---------------------------------------------------------------------------
$ws.PageSetup.Orientation = 1
$ws.Range("I3").Select()
